$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "DamagePerAttack"
$ws.Range("C1").Value = "TimePerAttack"
$ws.Range("D1").Value = "AttackRange"
$ws.Range("E2").Value = "Str/DamagePerAttack:+10;Dex/TimePerAttack:-1;"
